$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Target edit (per the commit diff):
#   "...Anul I la Cibernetică Economic –Nota de intrare 9,49."
#     -> "...Anul I la Informatica Economica –Nota de intrare 9,49."
#
# Underlying runs (paragraph "2019 –Academia de Studii Economice-Facultatea
# de Cibernetică, Statistică și Informatică Economică  Anul I la
# Cibernetică Economic –Nota de intrare 9,49.") around the edit site:
#   run6 " Anul I la "      (unchanged)
#   run7 "Cibernetic"       -> "Informatica Economica"
#   run8 "ă"                -> deleted
#   run9 " Economic –"      -> " –"
#   run10 "Nota"             (unchanged)
#   run11 " de intrare 9,49." (unchanged)
#
# run6, run7, run9, run10, run11 (and run0/run1 earlier in the paragraph)
# all share identical run formatting (color 000000 / lang ro). This
# engine recombines same-formatted adjacent runs across the whole
# paragraph whenever any run's text in that paragraph is edited, so we
# temporarily give each of those same-format runs a distinct dummy
# font color before editing (preventing unwanted merges), then restore
# the original color afterwards -- restoring a Font.Color alone does not
# trigger the recombination check.
# ---------------------------------------------------------------------

$findRange = $d.Content
$found = $findRange.Find.Execute("Cibernetică Economic –", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "target phrase not found"
}
$base = $findRange.Start

# Sanity-check the run layout we rely on.
if ($d.Range($base, $base + 10).Text -ne "Cibernetic") { throw "run7 mismatch" }
if ($d.Range($base + 10, $base + 11).Text -ne "ă") { throw "run8 mismatch" }
if ($d.Range($base + 11, $base + 22).Text -ne " Economic –") { throw "run9 mismatch" }
if ($d.Range($base + 22, $base + 26).Text -ne "Nota") { throw "run10 mismatch" }
if ($d.Range($base - 11, $base).Text -ne " Anul I la ") { throw "run6 mismatch" }
if ($d.Range($base + 26, $base + 43).Text -ne " de intrare 9,49.") { throw "run11 mismatch" }

# --- Phase 1: break apart every same-format run pair we must not let
#     the engine recombine while we edit this paragraph. ---
$dummy1 = 1
$dummy2 = 2
$dummy3 = 3
$dummy4 = 4
$dummy5 = 5

$rRun6  = $d.Range($base - 11, $base)
$rRun7  = $d.Range($base, $base + 10)
$rRun9  = $d.Range($base + 11, $base + 22)
$rRun10 = $d.Range($base + 22, $base + 26)
$rRun11 = $d.Range($base + 26, $base + 43)

$rRun6.Font.Color  = $dummy1
$rRun7.Font.Color  = $dummy2
$rRun9.Font.Color  = $dummy3
$rRun10.Font.Color = $dummy4
$rRun11.Font.Color = $dummy5

# --- Phase 2: perform the actual content edits. ---

# run9 " Economic –" -> " –"
$rRun9b = $d.Range($base + 11, $base + 22)
$rRun9b.Text = " –"

# run8 "ă" -> delete
$rRun8 = $d.Range($base + 10, $base + 11)
$rRun8.Delete()

# run7 "Cibernetic" -> "Informatica Economica"
$rRun7b = $d.Range($base, $base + 10)
$rRun7b.Text = "Informatica Economica"

# --- Phase 3: restore original formatting (color 000000) on every run
#     we perturbed, now re-resolved to their (shifted) final ranges. ---
$newRun7End = $base + ("Informatica Economica").Length
$newRun9End = $newRun7End + (" –").Length

$rRun6f  = $d.Range($base - 11, $base)
$rRun7f  = $d.Range($base, $newRun7End)
$rRun9f  = $d.Range($newRun7End, $newRun9End)
$rRun10f = $d.Range($newRun9End, $newRun9End + 4)
$rRun11f = $d.Range($newRun9End + 4, $newRun9End + 4 + (" de intrare 9,49.").Length)

$rRun6f.Font.Color  = 0
$rRun7f.Font.Color  = 0
$rRun9f.Font.Color  = 0
$rRun10f.Font.Color = 0
$rRun11f.Font.Color = 0

Write-Host "final: [$($d.Range($base - 11, $newRun9End + 4 + (" de intrare 9,49.").Length).Text)]"
